$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

$ws.Range("A$row").Value = "D3G7WT"
$ws.Range("B$row").Value = "Engranaje de combinación para impresora Epson"
$ws.Range("C$row").Value = "LX 300+"
$ws.Range("D$row").Value = 0
$ws.Range("E$row").Value = 50000
$ws.Range("F$row").Value = 9
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Formula = "=(E$row-D$row)*G$row"
$ws.Range("I$row").Formula = "=D$row*F$row"
$ws.Range("J$row").Value = 0
